$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 2005.6666
$ws.Range("H86").Value = 93751870
$ws.Range("I86").Value = 214287090
$ws.Range("J86").Value = 2266
$ws.Range("K86").Value = 214287090
$ws.Range("L86").Value = 2266
$ws.Range("M86").Value = -214285967
$ws.Range("N86").Value = -4512
$ws.Range("H89").Value = 93751870
$ws.Range("I89").Value = 214287090
$ws.Range("J89").Value = 2266
$ws.Range("K89").Value = 1071435450
$ws.Range("L89").Value = 11330
$ws.Range("M89").Value = -1071429834
$ws.Range("N89").Value = -22562
$ws.Range("H95").Value = 48888.8
$ws.Range("J95").Value = 48888.8
$ws.Range("L95").Value = 48888.8
$ws.Range("N95").Value = -54380.8
$ws.Range("H106").Value = 120001800
$ws.Range("I106").Value = 150001250
$ws.Range("K106").Value = 150001250
$ws.Range("M106").Value = -150000619
$ws.Range("H138").Value = 2363.2559
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2363.2559
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 7089.7677
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -17369.7677

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2318.5334
$ws.Range("I61").Value = 2059.8462
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2059.8462
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -1847.8462
$ws.Range("N61").Value = -4424
$ws.Range("I102").Value = 3123.75
$ws.Range("K102").Value = 3123.75
$ws.Range("M102").Value = -1501.75
$ws.Range("H128").Value = 35289.5
$ws.Range("J128").Value = 35289.5
$ws.Range("L128").Value = 35289.5
$ws.Range("N128").Value = -45249.5
$ws.Range("H132").Value = 3079.2222
$ws.Range("I132").Value = 1607.7273
$ws.Range("J132").Value = 5391.5713
$ws.Range("K132").Value = 4823.1819
$ws.Range("L132").Value = 16174.7139
$ws.Range("M132").Value = -2293.1819
$ws.Range("N132").Value = -21234.7139
$ws.Range("H136").Value = 2318.5334
$ws.Range("I136").Value = 2059.8462
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 6179.5386
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -3629.5386
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5437048
$ws.Range("I105").Value = 5684123
$ws.Range("J105").Value = 1400
$ws.Range("K105").Value = 5684123
$ws.Range("L105").Value = 1400
$ws.Range("M105").Value = -5682376
$ws.Range("N105").Value = -4894
$ws.Range("H131").Value = 40224
$ws.Range("J131").Value = 40224
$ws.Range("L131").Value = 40224
$ws.Range("N131").Value = -50304
$ws.Range("H133").Value = 39966.938
$ws.Range("J133").Value = 39966.938
$ws.Range("L133").Value = 39966.938
$ws.Range("N133").Value = -50086.938

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5691.807
$ws.Range("I31").Value = 1207.6857
$ws.Range("K31").Value = 1207.6857
$ws.Range("M31").Value = -912.6857
$ws.Range("H34").Value = 5691.807
$ws.Range("I34").Value = 1207.6857
$ws.Range("K34").Value = 1207.6857
$ws.Range("M34").Value = -1005.6857
$ws.Range("H58").Value = 1355.25
$ws.Range("I58").Value = 1361.2354
$ws.Range("J58").Value = 1321.3334
$ws.Range("K58").Value = 1361.2354
$ws.Range("L58").Value = 1321.3334
$ws.Range("M58").Value = -1158.2354
$ws.Range("N58").Value = -1727.3334
$ws.Range("H105").Value = 1940
$ws.Range("I105").Value = 1700
$ws.Range("J105").Value = 2260
$ws.Range("K105").Value = 1700
$ws.Range("L105").Value = 2260
$ws.Range("M105").Value = 47
$ws.Range("N105").Value = -5754
$ws.Range("H107").Value = 4167618.8
$ws.Range("I107").Value = 6250728
$ws.Range("J107").Value = 1400
$ws.Range("K107").Value = 6250728
$ws.Range("L107").Value = 1400
$ws.Range("M107").Value = -6248808
$ws.Range("N107").Value = -5240
$ws.Range("H136").Value = 1355.25
$ws.Range("I136").Value = 1361.2354
$ws.Range("J136").Value = 1321.3334
$ws.Range("K136").Value = 4083.7062
$ws.Range("L136").Value = 3964.0002
$ws.Range("M136").Value = -1533.7062
$ws.Range("N136").Value = -9064.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 23336334
$ws.Range("H31").Value = 1165
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 1165
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3495
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -4071
$ws.Range("H110").Value = 14124.706
$ws.Range("J110").Value = 14124.706
$ws.Range("L110").Value = 42374.118
$ws.Range("N110").Value = -50554.118
$ws.Range("H136").Value = 2066.5
$ws.Range("I136").Value = 1099.75
$ws.Range("K136").Value = 3299.25
$ws.Range("M136").Value = 1800.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1667450.5
$ws.Range("I3").Value = 5000250
$ws.Range("J3").Value = 1050.75
$ws.Range("K3").Value = 5000250
$ws.Range("L3").Value = 1050.75
$ws.Range("M3").Value = -5000134
$ws.Range("N3").Value = -1282.75
$ws.Range("H10").Value = 8611287
$ws.Range("I10").Value = 20000000
$ws.Range("J10").Value = 69752
$ws.Range("K10").Value = 20000000
$ws.Range("L10").Value = 69752
$ws.Range("M10").Value = -19999831
$ws.Range("N10").Value = -70090
$ws.Range("H14").Value = 18334502
$ws.Range("I14").Value = 18334502
$ws.Range("K14").Value = 18334502
$ws.Range("M14").Value = -18334334
$ws.Range("H21").Value = 5775
$ws.Range("J21").Value = 5775
$ws.Range("L21").Value = 5775
$ws.Range("N21").Value = -6121
$ws.Range("H30").Value = 5775
$ws.Range("J30").Value = 5775
$ws.Range("L30").Value = 5775
$ws.Range("N30").Value = -5985
$ws.Range("H112").Value = 55999.25
$ws.Range("J112").Value = 55999.25
$ws.Range("L112").Value = 55999.25
$ws.Range("N112").Value = -58215.25
$ws.Range("H114").Value = 34044
$ws.Range("J114").Value = 34044
$ws.Range("L114").Value = 34044
$ws.Range("N114").Value = -42722
$ws.Range("H116").Value = 38999.25
$ws.Range("J116").Value = 38999.25
$ws.Range("L116").Value = 38999.25
$ws.Range("N116").Value = -48177.25
$ws.Range("H117").Value = 64873.332
$ws.Range("J117").Value = 64873.332
$ws.Range("L117").Value = 64873.332
$ws.Range("N117").Value = -71757.33199999999
$ws.Range("H118").Value = 74810
$ws.Range("J118").Value = 74810
$ws.Range("L118").Value = 74810
$ws.Range("N118").Value = -78124

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 13891232
$ws.Range("I136").Value = 2284
$ws.Range("J136").Value = 27780180
$ws.Range("K136").Value = 6852
$ws.Range("L136").Value = 83340540
$ws.Range("M136").Value = -4302
$ws.Range("N136").Value = -83345640

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 145736000
$ws.Range("J2").Value = 200030400
$ws.Range("L2").Value = 200030400
$ws.Range("N2").Value = -200030624
$ws.Range("H3").Value = 37502
$ws.Range("I3").Value = 5500
$ws.Range("J3").Value = 69504
$ws.Range("K3").Value = 5500
$ws.Range("L3").Value = 69504
$ws.Range("M3").Value = -5386
$ws.Range("N3").Value = -69732
$ws.Range("H4").Value = 69668.664
$ws.Range("J4").Value = 69668.664
$ws.Range("L4").Value = 69668.664
$ws.Range("N4").Value = -69894.664
$ws.Range("H10").Value = 57603.6
$ws.Range("I10").Value = 8000
$ws.Range("J10").Value = 70004.5
$ws.Range("K10").Value = 8000
$ws.Range("L10").Value = 70004.5
$ws.Range("M10").Value = -7831
$ws.Range("N10").Value = -70342.5
$ws.Range("H14").Value = 4919800.5
$ws.Range("I14").Value = 70000
$ws.Range("J14").Value = 5458667
$ws.Range("K14").Value = 70000
$ws.Range("L14").Value = 5458667
$ws.Range("M14").Value = -69832
$ws.Range("N14").Value = -5459003
$ws.Range("H21").Value = 70017
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 70017
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 70017
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -70487
$ws.Range("H32").Value = 500004000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 500004000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 500004000
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -500004634
$ws.Range("H33").Value = 70021
$ws.Range("J33").Value = 70021
$ws.Range("L33").Value = 70021
$ws.Range("N33").Value = -70521
$ws.Range("H35").Value = 70017
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 70017
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 70017
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -70597
$ws.Range("H36").Value = 70021
$ws.Range("J36").Value = 70021
$ws.Range("L36").Value = 70021
$ws.Range("N36").Value = -70521
$ws.Range("H39").Value = 70048.5
$ws.Range("I39").Value = 70048
$ws.Range("K39").Value = 70048
$ws.Range("M39").Value = -69635
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H81").Value = 5476.7
$ws.Range("I81").Value = 4824.2856
$ws.Range("J81").Value = 6999
$ws.Range("K81").Value = 9648.5712
$ws.Range("L81").Value = 13998
$ws.Range("M81").Value = -8587.5712
$ws.Range("N81").Value = -16120
$ws.Range("H84").Value = 5476.7
$ws.Range("I84").Value = 4824.2856
$ws.Range("J84").Value = 6999
$ws.Range("K84").Value = 48242.856
$ws.Range("L84").Value = 69990
$ws.Range("M84").Value = -42938.856
$ws.Range("N84").Value = -80598
$ws.Range("H101").Value = 15958.25
$ws.Range("J101").Value = 15958.25
$ws.Range("L101").Value = 15958.25
$ws.Range("N101").Value = -22448.25
